$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 08:50"

# Row 17
$ws.Range("B17").Value = 9634
$ws.Range("C17").Value = 16
$ws.Range("E17").Value = 8890

# Row 30
$ws.Range("A30").Value = "Polonia"
$ws.Range("B30").Value = 2132
$ws.Range("C30").Value = 77
$ws.Range("D30").Value = 7
$ws.Range("E30").Value = 2094
$ws.Range("F30").Value = 3
$ws.Range("H30").Value = 31

# Row 31
$ws.Range("A31").Value = "Rumania"
$ws.Range("B31").Value = 2109
$ws.Range("D31").Value = 209
$ws.Range("E31").Value = 1835
$ws.Range("F31").Value = 33
$ws.Range("H31").Value = 65

# Row 37
$ws.Range("A37").Value = "Tailandia"
$ws.Range("B37").Value = 1651
$ws.Range("C37").Value = 127
$ws.Range("D37").Value = 229
$ws.Range("E37").Value = 1412
$ws.Range("F37").Value = 11
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = 10

# Row 38
$ws.Range("A38").Value = "Filipinas"
$ws.Range("B38").Value = 1546
$ws.Range("D38").Value = 42
$ws.Range("E38").Value = 1426
$ws.Range("F38").Value = 1
$ws.Range("H38").Value = 78

# Row 73
$ws.Range("A73").Value = "Bosnia y Herzegovina"
$ws.Range("B73").Value = 379
$ws.Range("C73").Value = 11
$ws.Range("D73").Value = 17
$ws.Range("E73").Value = 352
$ws.Range("F73").Value = 1
$ws.Range("H73").Value = 10

# Row 74
$ws.Range("A74").Value = "Letonia"
$ws.Range("B74").Value = 376
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 1
$ws.Range("E74").Value = 375
$ws.Range("F74").Value = 3
$ws.Range("H74").Value = 0

# Row 80
$ws.Range("A80").Value = "Taiwan"
$ws.Range("B80").Value = 322
$ws.Range("C80").Value = 16
$ws.Range("D80").Value = 39
$ws.Range("E80").Value = 278
$ws.Range("F80").Value = 0
$ws.Range("H80").Value = 5

# Row 81
$ws.Range("A81").Value = "Uruguay"
$ws.Range("B81").Value = 320
$ws.Range("D81").Value = 25
$ws.Range("E81").Value = 294
$ws.Range("F81").Value = 9
$ws.Range("H81").Value = 1

# Row 93
$ws.Range("A93").Value = "Camerun"
$ws.Range("B93").Value = 193
$ws.Range("C93").Value = 54
$ws.Range("D93").Value = 5
$ws.Range("E93").Value = 182
$ws.Range("F93").Value = 0
$ws.Range("H93").Value = 6

# Row 94
$ws.Range("A94").Value = "Oman"
$ws.Range("B94").Value = 179
$ws.Range("D94").Value = 29
$ws.Range("E94").Value = 150
$ws.Range("F94").Value = 3
$ws.Range("H94").Value = 0

# Row 95
$ws.Range("A95").Value = "Afganistan"
$ws.Range("D95").Value = 2
$ws.Range("E95").Value = 164
$ws.Range("F95").Value = 0

# Row 96
$ws.Range("A96").Value = "Cuba"
$ws.Range("B96").Value = 170
$ws.Range("D96").Value = 4
$ws.Range("E96").Value = 162
$ws.Range("F96").Value = 2
$ws.Range("H96").Value = 4

# Row 97
$ws.Range("B97").Value = 169
$ws.Range("C97").Value = 1
$ws.Range("D97").Value = 74
$ws.Range("E97").Value = 95
$ws.Range("F97").Value = 3

# Row 98
$ws.Range("A98").Value = "Costa de Marfil"
$ws.Range("B98").Value = 168
$ws.Range("D98").Value = 6
$ws.Range("E98").Value = 161
$ws.Range("H98").Value = 1

# Row 99
$ws.Range("A99").Value = "Senegal"
$ws.Range("B99").Value = 162
$ws.Range("D99").Value = 27
$ws.Range("E99").Value = 135
$ws.Range("F99").Value = 0

# Row 100
$ws.Range("A100").Value = "Malta"
$ws.Range("B100").Value = 156
$ws.Range("E100").Value = 154
$ws.Range("F100").Value = 4
$ws.Range("H100").Value = 0

# Row 101
$ws.Range("A101").Value = "Ghana"
$ws.Range("D101").Value = 2
$ws.Range("E101").Value = 145
$ws.Range("F101").Value = 1
$ws.Range("H101").Value = 5

# Row 102
$ws.Range("A102").Value = "Bielorrusia"
$ws.Range("B102").Value = 152
$ws.Range("C102").Value = 0
$ws.Range("D102").Value = 32
$ws.Range("E102").Value = 120
$ws.Range("F102").Value = 2
$ws.Range("H102").Value = 0

# Row 103
$ws.Range("A103").Value = "Uzbekistan"
$ws.Range("B103").Value = 150
$ws.Range("C103").Value = 1
$ws.Range("D103").Value = 7
$ws.Range("E103").Value = 141
$ws.Range("F103").Value = 8
$ws.Range("H103").Value = 2

# Row 104
$ws.Range("A104").Value = "Honduras"
$ws.Range("B104").Value = 141
$ws.Range("C104").Value = 2
$ws.Range("D104").Value = 3
$ws.Range("E104").Value = 131
$ws.Range("F104").Value = 4
$ws.Range("H104").Value = 7

# Row 111
$ws.Range("A111").Value = "Georgia"
$ws.Range("B111").Value = 108
$ws.Range("C111").Value = 5
$ws.Range("D111").Value = 21
$ws.Range("E111").Value = 87
$ws.Range("F111").Value = 6

# Row 112
$ws.Range("A112").Value = "Kirguistan"
$ws.Range("C112").Value = 13
$ws.Range("D112").Value = 3
$ws.Range("E112").Value = 104
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 0

# Row 113
$ws.Range("A113").Value = "Bolivia"
$ws.Range("C113").Value = 10
$ws.Range("D113").Value = 0
$ws.Range("E113").Value = 101
$ws.Range("F113").Value = 3
$ws.Range("G113").Value = 2
$ws.Range("H113").Value = 6

# Row 114
$ws.Range("A114").Value = "Camboya"
$ws.Range("B114").Value = 107
$ws.Range("D114").Value = 21
$ws.Range("E114").Value = 86
$ws.Range("F114").Value = 1
$ws.Range("H114").Value = 0

# Row 115
$ws.Range("A115").Value = "Guadalupe"
$ws.Range("B115").Value = 106
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 17
$ws.Range("E115").Value = 85
$ws.Range("F115").Value = 10
$ws.Range("H115").Value = 4

# Row 116
$ws.Range("A116").Value = "Montenegro"
$ws.Range("B116").Value = 105
$ws.Range("C116").Value = 14
$ws.Range("D116").Value = 0
$ws.Range("E116").Value = 104
$ws.Range("F116").Value = 1
$ws.Range("H116").Value = 1

# Row 130
$ws.Range("B130").Value = 46
$ws.Range("C130").Value = 3
$ws.Range("E130").Value = 46
$ws.Range("F130").Value = 6

# Row 154
$ws.Range("A154").Value = "Eritrea"

# Row 155
$ws.Range("A155").Value = "Nueva Caledonia"

# Row 158
$ws.Range("A158").Value = "Bahamas"
$ws.Range("D158").Value = 1
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 0

# Row 159
$ws.Range("A159").Value = "Birmania"
$ws.Range("D159").Value = 0
$ws.Range("G159").Value = 1
$ws.Range("H159").Value = 1

# Row 160
$ws.Range("A160").Value = "Guinea Ecuatorial"

# Row 161
$ws.Range("A161").Value = "Dominica"

# Row 167
$ws.Range("A167").Value = "Siria"
$ws.Range("D167").Value = 0
$ws.Range("H167").Value = 2

# Row 168
$ws.Range("A168").Value = "Groenlandia"
$ws.Range("D168").Value = 2
$ws.Range("H168").Value = 0

# Row 170
$ws.Range("A170").Value = "Laos"
$ws.Range("C170").Value = 1

# Row 171
$ws.Range("A171").Value = "Granada"
$ws.Range("C171").Value = 0

# Row 173
$ws.Range("A173").Value = "Surinam"

# Row 174
$ws.Range("A174").Value = "Libia"

# Row 175
$ws.Range("A175").Value = "Guinea-Bisau"

# Row 176
$ws.Range("A176").Value = "Mozambique"

# Row 182
$ws.Range("A182").Value = "Santa Sede"

# Row 183
$ws.Range("A183").Value = "San Martin (Parte Holandesa)"

# Row 198
$ws.Range("A198").Value = "Botsuana"

# Row 199
$ws.Range("A199").Value = "Belice"

# Row 200
$ws.Range("A200").Value = "Liberia"
$ws.Range("C200").Value = 0

# Row 201
$ws.Range("A201").Value = "Islas Virgenes Britanicas"
$ws.Range("C201").Value = 1
